$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix punctuation in proveedor names (comma -> period) ---
$ws.Range("E24").Value = 'IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA'
$ws.Range("F24").Value = 'IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA'
$ws.Range("E44").Value = 'IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA'
$ws.Range("F44").Value = 'IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA'
$ws.Range("E53").Value = 'FERNANDEZ MARIO H. GALLICET OSCAR M'
$ws.Range("E54").Value = 'IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA'
$ws.Range("F54").Value = 'IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA'
$ws.Range("E55").Value = 'MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO'
$ws.Range("E92").Value = 'FERNANDEZ MARIO H. GALLICET OSCAR M'
$ws.Range("E93").Value = 'MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO'

# --- Fix floating point formatting in Importe column (H) ---
# Force text number format so Excel keeps these as text strings, matching
# the original data (avoids locale-driven reinterpretation as numbers).
$ws.Range("H2:H128").NumberFormat = "@"
$ws.Range("H2").Value = "99.63"
$ws.Range("H3").Value = "58000.00"
$ws.Range("H4").Value = "1292.40"
$ws.Range("H5").Value = "4986.59"
$ws.Range("H6").Value = "96723.70"
$ws.Range("H7").Value = "117655.48"
$ws.Range("H8").Value = "11227.50"
$ws.Range("H9").Value = "285.00"
$ws.Range("H10").Value = "187.96"
$ws.Range("H11").Value = "8255.95"
$ws.Range("H12").Value = "319.00"
$ws.Range("H13").Value = "12432.56"
$ws.Range("H14").Value = "29430.00"
$ws.Range("H15").Value = "2000.00"
$ws.Range("H16").Value = "320.00"
$ws.Range("H17").Value = "330.00"
$ws.Range("H18").Value = "961.00"
$ws.Range("H19").Value = "410.00"
$ws.Range("H20").Value = "104.00"
$ws.Range("H21").Value = "875.00"
$ws.Range("H22").Value = "5193.69"
$ws.Range("H23").Value = "265.86"
$ws.Range("H24").Value = "7.89"
$ws.Range("H25").Value = "78.00"
$ws.Range("H26").Value = "590.00"
$ws.Range("H27").Value = "48298.35"
$ws.Range("H28").Value = "58.50"
$ws.Range("H29").Value = "3329.00"
$ws.Range("H30").Value = "4240.00"
$ws.Range("H31").Value = "257.88"
$ws.Range("H32").Value = "3850.00"
$ws.Range("H33").Value = "39.68"
$ws.Range("H34").Value = "6311.89"
$ws.Range("H35").Value = "2375.00"
$ws.Range("H36").Value = "304.00"
$ws.Range("H37").Value = "196.00"
$ws.Range("H38").Value = "444.00"
$ws.Range("H39").Value = "7020.00"
$ws.Range("H40").Value = "1757.00"
$ws.Range("H41").Value = "48.16"
$ws.Range("H42").Value = "38434.00"
$ws.Range("H43").Value = "43.50"
$ws.Range("H44").Value = "60.46"
$ws.Range("H45").Value = "1680.00"
$ws.Range("H46").Value = "1600.00"
$ws.Range("H47").Value = "284.10"
$ws.Range("H48").Value = "27180.00"
$ws.Range("H49").Value = "550.00"
$ws.Range("H50").Value = "2370.00"
$ws.Range("H51").Value = "250.00"
$ws.Range("H52").Value = "1441.00"
$ws.Range("H53").Value = "950.00"
$ws.Range("H54").Value = "129.26"
$ws.Range("H55").Value = "420.00"
$ws.Range("H56").Value = "6065.00"
$ws.Range("H57").Value = "2480.00"
$ws.Range("H58").Value = "92500.00"
$ws.Range("H59").Value = "132870.00"
$ws.Range("H60").Value = "5751.60"
$ws.Range("H61").Value = "8030.00"
$ws.Range("H62").Value = "97.93"
$ws.Range("H63").Value = "11545.20"
$ws.Range("H64").Value = "2688.00"
$ws.Range("H65").Value = "617.60"
$ws.Range("H66").Value = "1476.00"
$ws.Range("H67").Value = "18.42"
$ws.Range("H68").Value = "4276.80"
$ws.Range("H69").Value = "1500.00"
$ws.Range("H70").Value = "1508.22"
$ws.Range("H71").Value = "198.00"
$ws.Range("H72").Value = "259.52"
$ws.Range("H73").Value = "450.00"
$ws.Range("H74").Value = "25885.00"
$ws.Range("H75").Value = "4042.74"
$ws.Range("H76").Value = "759500.00"
$ws.Range("H77").Value = "160.00"
$ws.Range("H78").Value = "14.00"
$ws.Range("H79").Value = "1304.86"
$ws.Range("H80").Value = "7950.00"
$ws.Range("H81").Value = "200.00"
$ws.Range("H82").Value = "500.00"
$ws.Range("H83").Value = "445.00"
$ws.Range("H84").Value = "1000.00"
$ws.Range("H85").Value = "500.00"
$ws.Range("H86").Value = "950.00"
$ws.Range("H87").Value = "500.00"
$ws.Range("H88").Value = "4020.00"
$ws.Range("H89").Value = "2250.00"
$ws.Range("H90").Value = "200.00"
$ws.Range("H91").Value = "440.00"
$ws.Range("H92").Value = "340.00"
$ws.Range("H93").Value = "425.00"
$ws.Range("H94").Value = "10792.71"
$ws.Range("H95").Value = "7300.00"
$ws.Range("H96").Value = "81.08"
$ws.Range("H97").Value = "1390.41"
$ws.Range("H98").Value = "7770.00"
$ws.Range("H99").Value = "21.00"
$ws.Range("H100").Value = "2376.00"
$ws.Range("H101").Value = "200.00"
$ws.Range("H102").Value = "549.00"
$ws.Range("H103").Value = "509.00"
$ws.Range("H104").Value = "284.86"
$ws.Range("H105").Value = "190.00"
$ws.Range("H106").Value = "1557.50"
$ws.Range("H107").Value = "194.34"
$ws.Range("H108").Value = "645.00"
$ws.Range("H109").Value = "231.00"
$ws.Range("H110").Value = "325.77"
$ws.Range("H111").Value = "5288.09"
$ws.Range("H112").Value = "1680.69"
$ws.Range("H113").Value = "21910.00"
$ws.Range("H114").Value = "315.00"
$ws.Range("H115").Value = "44500.00"
$ws.Range("H116").Value = "325271.49"
$ws.Range("H117").Value = "1124.09"
$ws.Range("H118").Value = "402800.00"
$ws.Range("H119").Value = "289456.00"
$ws.Range("H120").Value = "381222.00"
$ws.Range("H121").Value = "377968.00"
$ws.Range("H122").Value = "414000.00"
$ws.Range("H123").Value = "409048.00"
$ws.Range("H124").Value = "380000.00"
$ws.Range("H125").Value = "100.00"
$ws.Range("H126").Value = "828797.44"
$ws.Range("H127").Value = "16900.00"
$ws.Range("H128").Value = "3302.40"
